$d = $word.ActiveDocument
$apos = [char]0x2019

# --- 1. Split "see a brief description. Google's..." run into three runs,
#        inserting new sentence and moving the _GoBack bookmark in between. ---

# Remove the bookmark from its old location (end of document) first so that
# Find operations below are not affected by it, then we re-add it later.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

$target = "see a brief description. Google" + $apos + "s natural language processing RESTful API is used to perform an entity analysis on the article description to produce Wikipedia links of the relevant topics."
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute($target, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target sentence"
}

$sentenceStart = $rng.Start
# "see a brief description. " occupies the first 26 characters of the found range.
$prefixLen = ("see a brief description. ").Length
$insertPoint = $sentenceStart + $prefixLen

$insertRng = $d.Range($insertPoint, $insertPoint)
$insertRng.InsertAfter("To complete the extra credit portion of the assignment, ")

# Re-create the bookmark right after the newly inserted text (collapsed range).
$bmStart = $insertPoint + ("To complete the extra credit portion of the assignment, ").Length
$bmRng = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmRng)

# --- 2. Merge "Top left: View " + "hen the application is first opened" ---
$d.Content.Find.Execute("Top left: View hen the application is first opened", $true, $false, $false, $false, $false, $true, 1, $false, "Top left: View hen the application is first opened", 2)

# --- 3. Merge "Author: Christopher Pearce" + ", " + "Robert " (first occurrence, before Farinelli/AsyncResponse) ---
$d.Content.Find.Execute("Author: Christopher Pearce, Robert Farinelli", $true, $false, $false, $false, $false, $true, 1, $false, "Author: Christopher Pearce, Robert Farinelli", 2)

# --- 4 & 5. Merge ", " + "Christopher Pearce" (two occurrences after "Author: Robert Farinelli") ---
$searchText = "Author: Robert Farinelli, Christopher Pearce"
$d.Content.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $searchText, 2)

# --- 6. Merge "Author: Christopher Pearce, " + "Robert " ---
$d.Content.Find.Execute("Author: Christopher Pearce, Robert Farinelli", $true, $false, $false, $false, $false, $true, 1, $false, "Author: Christopher Pearce, Robert Farinelli", 2)

# --- 7. Merge "Author: " + "Robert " ---
$d.Content.Find.Execute("Author: Robert Farinelli", $true, $false, $false, $false, $false, $true, 1, $false, "Author: Robert Farinelli", 2)
